$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 40: SRH vs RR - fill in player scores for row 52
$ws.Range("E52").Value = 60
$ws.Range("H52").Value = 100
$ws.Range("K52").Value = 30
$ws.Range("N52").Value = 80
$ws.Range("Q52").Value = 40
$ws.Range("T52").Value = 20
$ws.Range("W52").Value = 50
$ws.Range("Z52").Value = 0
$ws.Range("AC52").Value = 70

$excel.CalculateFullRebuild()
